$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: football-player factor description / short label
$ws.Range("A3").Value = 'Here is the factor description in terms of the features of a football player: ```The factor is very strongly positively associated with the feature that final third passes, very strongly positively associated with the feature that final third receptions, very strongly positively associated with the feature that smart passes, very strongly positively associated with the feature that adjusted, strongly positively associated with the feature that ground duels. The factor is weakly negatively associated with the feature that air duels, very weakly negatively associated with the feature that non-penalty expected goals. ```'
$ws.Range("B3").Value = 'creative passer vs aerial dueler'

# Row 4: country factor description / short label
$ws.Range("A4").Value = 'The existing names are: religious vs secular, confident vs doubtful. In this case, it is important that the name you now make is different from these names.
Here is the factor description in terms of the features of a country: ```The factor is very strongly positively associated with the feature that Are you an active member, inactive member, or not a member of a sport or recreational organization?, very strongly positively associated with the feature that Are you an active member, inactive member, or not a member of another type of organization?, strongly positively associated with the feature that Are you an active member, inactive member, or not a member of an art, music, or educational organization?, strongly positively associated with the feature that Are you an active member, inactive member, or not a member of a consumer organization?, strongly positively associated with the feature that Are you an active member, inactive member, or not a member of a labor union?. The factor is very strongly negatively associated with the feature that How frequently does alcohol consumption in the streets occur in your neighborhood?, very strongly negatively associated with the feature that Please tell me for the following statement how essential you think it is as a characteristic of democracy: Governments tax the rich and subsidize the poor., very strongly negatively associated with the feature that Have you been the victim of a crime during the past year?, strongly negatively associated with the feature that Has someone in your immediate family been the victim of a crime during the past year?, strongly negatively associated with the feature that Do you agree that science and technology will provide more opportunities for the next generation?. ```'
$ws.Range("B4").Value = 'community membership vs insecurity and crime'

# Give the new rows the same look (font/fill/border/wrap) as the existing
# example row above them (row 2), then size the rows to fit the new text.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows("3:3").RowHeight = 143.8
$ws.Rows("4:4").RowHeight = 352.7
